# Update feature importance values (decision tree model) per latest training run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.04142577817499517
$ws.Cells.Item(2, 3).Value = 0.009119430206315867
$ws.Cells.Item(2, 4).Value = 0.06241182299095948
$ws.Cells.Item(2, 5).Value = 0.03454793402522253
$ws.Cells.Item(2, 6).Value = 0.006746983357423617
$ws.Cells.Item(2, 7).Value = 0.01708981157364537
$ws.Cells.Item(2, 8).Value = 0.002754214343031071
$ws.Cells.Item(2, 9).Value = 0.007744415193710784
$ws.Cells.Item(2, 10).Value = 0.01122666382114934
$ws.Cells.Item(2, 11).Value = 0.0293444933506957
$ws.Cells.Item(2, 12).Value = 0.01874273300411592
$ws.Cells.Item(2, 13).Value = 0.1176529492183252
$ws.Cells.Item(2, 14).Value = 0.03598397767625214
$ws.Cells.Item(2, 15).Value = 0.03745775885144884
$ws.Cells.Item(2, 16).Value = 0.04585725763342981
$ws.Cells.Item(2, 17).Value = 0.001903016537504592
$ws.Cells.Item(2, 18).Value = 0.01602069595179173
$ws.Cells.Item(2, 19).Value = 0.01585599591014115
$ws.Cells.Item(2, 20).Value = 0.003552456232356023
$ws.Cells.Item(2, 21).Value = 0.001103554193260878
$ws.Cells.Item(2, 22).Value = 0.003244204487283796
$ws.Cells.Item(2, 23).Value = 0.02028403547405648
$ws.Cells.Item(2, 24).Value = 0.009219216014908488
$ws.Cells.Item(2, 25).Value = 0.01297439705326808
$ws.Cells.Item(2, 26).Value = 0.3789218855300409
$ws.Cells.Item(2, 27).Value = 0.006435430980247546
$ws.Cells.Item(2, 28).Value = 0.003940864234113369
$ws.Cells.Item(2, 29).Value = 0.0005986726383221138
$ws.Cells.Item(2, 30).Value = 0.03239904666028151
$ws.Cells.Item(2, 31).Value = 0.009558086861860977
$ws.Cells.Item(2, 32).Value = 0.005882217819841569

$ws.Cells.Item(3, 2).Value = 0.005430693145316895
$ws.Cells.Item(3, 3).Value = 0.01733014739952
$ws.Cells.Item(3, 4).Value = 0.01460380058829609
$ws.Cells.Item(3, 5).Value = 0.04127545179108195
$ws.Cells.Item(3, 6).Value = 0.00761940414211925
$ws.Cells.Item(3, 7).Value = 0.02690181287263051
$ws.Cells.Item(3, 8).Value = 0.01077276703173507
$ws.Cells.Item(3, 9).Value = 0.08821409945971646
$ws.Cells.Item(3, 10).Value = 0.06012404394019871
$ws.Cells.Item(3, 11).Value = 0.009147714987890382
$ws.Cells.Item(3, 12).Value = 0.05213656720489655
$ws.Cells.Item(3, 13).Value = 0.01658770328325978
$ws.Cells.Item(3, 14).Value = 0.01978779799041725
$ws.Cells.Item(3, 15).Value = 0.01666765660451682
$ws.Cells.Item(3, 16).Value = 0.01868114409032903
$ws.Cells.Item(3, 17).Value = 0.08745644202111151
$ws.Cells.Item(3, 18).Value = 0.002255771230602203
$ws.Cells.Item(3, 19).Value = 0.0223610516622046
$ws.Cells.Item(3, 20).Value = 0.0402277627222782
$ws.Cells.Item(3, 21).Value = 0.002155325804836443
$ws.Cells.Item(3, 22).Value = 0.009370791832040608
$ws.Cells.Item(3, 23).Value = 0.02797177703317714
$ws.Cells.Item(3, 24).Value = 0.01227857524741984
$ws.Cells.Item(3, 25).Value = 0.03226727897835125
$ws.Cells.Item(3, 26).Value = 0.01011943064000844
$ws.Cells.Item(3, 27).Value = 0.03125092536676131
$ws.Cells.Item(3, 28).Value = 0.08521175204645025
$ws.Cells.Item(3, 29).Value = 0.002615321105858613
$ws.Cells.Item(3, 30).Value = 0.1361504474633851
$ws.Cells.Item(3, 31).Value = 0.04333906054307166
$ws.Cells.Item(3, 32).Value = 0.04968748177051784

$ws.Cells.Item(4, 2).Value = 0.0375367542204979
$ws.Cells.Item(4, 3).Value = 0.01104799532427805
$ws.Cells.Item(4, 4).Value = 0.02124953338710933
$ws.Cells.Item(4, 5).Value = 0.1461380717821183
$ws.Cells.Item(4, 6).Value = 0.004379993015605014
$ws.Cells.Item(4, 7).Value = 0.01444427463941971
$ws.Cells.Item(4, 8).Value = 0.05008092761301727
$ws.Cells.Item(4, 9).Value = 0.00611441201114897
$ws.Cells.Item(4, 10).Value = 0.007412059891202965
$ws.Cells.Item(4, 11).Value = 0.01639113273721011
$ws.Cells.Item(4, 12).Value = 0.05538818647924072
$ws.Cells.Item(4, 13).Value = 0.05000926335668898
$ws.Cells.Item(4, 14).Value = 0.01235488903822852
$ws.Cells.Item(4, 15).Value = 0.003072695031767502
$ws.Cells.Item(4, 16).Value = 0.01278188093152356
$ws.Cells.Item(4, 17).Value = 0.05419485961855384
$ws.Cells.Item(4, 18).Value = 0.006561300411942304
$ws.Cells.Item(4, 19).Value = 0.001195422787738216
$ws.Cells.Item(4, 20).Value = 0.006201744968329247
$ws.Cells.Item(4, 21).Value = 0.003847427148078504
$ws.Cells.Item(4, 22).Value = 0.009697053459384644
$ws.Cells.Item(4, 23).Value = 0.02017929451578506
$ws.Cells.Item(4, 24).Value = 0.03399465607717327
$ws.Cells.Item(4, 25).Value = 0.003557322283427295
$ws.Cells.Item(4, 26).Value = 0.2060217991627378
$ws.Cells.Item(4, 27).Value = 0.04147274177975339
$ws.Cells.Item(4, 28).Value = 0.00405944574485724
$ws.Cells.Item(4, 29).Value = 0.009948158276379812
$ws.Cells.Item(4, 30).Value = 0.03818103056963332
$ws.Cells.Item(4, 31).Value = 0.03164705416915764
$ws.Cells.Item(4, 32).Value = 0.08083861956801182

$ws.Cells.Item(5, 2).Value = 0.003654232166171411
$ws.Cells.Item(5, 3).Value = 0.002158544405972808
$ws.Cells.Item(5, 4).Value = 0.01183967758953136
$ws.Cells.Item(5, 5).Value = 0.01866090832493984
$ws.Cells.Item(5, 6).Value = 0.01446280946140693
$ws.Cells.Item(5, 7).Value = 0.007229007585270534
$ws.Cells.Item(5, 8).Value = 0.0003810451720702955
$ws.Cells.Item(5, 9).Value = 0.01278701192038035
$ws.Cells.Item(5, 10).Value = 0.1476211048891866
$ws.Cells.Item(5, 11).Value = 0.1089982620909422
$ws.Cells.Item(5, 12).Value = 0.0005195480759701095
$ws.Cells.Item(5, 13).Value = 0.05046100159478484
$ws.Cells.Item(5, 14).Value = 0.01469944268240226
$ws.Cells.Item(5, 15).Value = 0.005796035707624606
$ws.Cells.Item(5, 16).Value = 0.01997210143578817
$ws.Cells.Item(5, 17).Value = 0.009595571741369078
$ws.Cells.Item(5, 18).Value = 0.02507843072808258
$ws.Cells.Item(5, 19).Value = 0.06044767317500175
$ws.Cells.Item(5, 20).Value = 0.006199195465047543
$ws.Cells.Item(5, 21).Value = 0.02335360735967333
$ws.Cells.Item(5, 22).Value = 0.01355941691592833
$ws.Cells.Item(5, 23).Value = 0.01721846991812657
$ws.Cells.Item(5, 24).Value = 0.2411835088666599
$ws.Cells.Item(5, 25).Value = 0.000063189005652573672420883877
$ws.Cells.Item(5, 26).Value = 0.03421501058572696
$ws.Cells.Item(5, 27).Value = 0.002544259671756123
$ws.Cells.Item(5, 28).Value = 0.002345196231393964
$ws.Cells.Item(5, 29).Value = 0.004727319463771576
$ws.Cells.Item(5, 30).Value = 0.06557196696475304
$ws.Cells.Item(5, 31).Value = 0.02104673627157618
$ws.Cells.Item(5, 32).Value = 0.05360971453303842

$ws.Cells.Item(6, 2).Value = 0.002160616185580105
$ws.Cells.Item(6, 3).Value = 0.02885305528827399
$ws.Cells.Item(6, 4).Value = 0.004749989247254319
$ws.Cells.Item(6, 5).Value = 0.003508283595324904
$ws.Cells.Item(6, 6).Value = 0.02134129292158206
$ws.Cells.Item(6, 7).Value = 0.009636866492675016
$ws.Cells.Item(6, 8).Value = 0.0000228886814280114416050773
$ws.Cells.Item(6, 9).Value = 0.01205193860758435
$ws.Cells.Item(6, 10).Value = 0.01067785276152942
$ws.Cells.Item(6, 11).Value = 0.001736406199645091
$ws.Cells.Item(6, 12).Value = 0.0522348654203845
$ws.Cells.Item(6, 13).Value = 0.02987110892152077
$ws.Cells.Item(6, 14).Value = 0.009475843717048736
$ws.Cells.Item(6, 15).Value = 0.008091326248291021
$ws.Cells.Item(6, 16).Value = 0.04669995423115453
$ws.Cells.Item(6, 17).Value = 0.1057454707948054
$ws.Cells.Item(6, 18).Value = 0.01219480080115458
$ws.Cells.Item(6, 19).Value = 0.06861092574414444
$ws.Cells.Item(6, 20).Value = 0.001841006657040083
$ws.Cells.Item(6, 21).Value = 0.006958404391208561
$ws.Cells.Item(6, 22).Value = 0.07486274267063055
$ws.Cells.Item(6, 23).Value = 0.01712233474334489
$ws.Cells.Item(6, 24).Value = 0.09800816490071609
$ws.Cells.Item(6, 25).Value = 0.002262070375608591
$ws.Cells.Item(6, 26).Value = 0.008534597867902777
$ws.Cells.Item(6, 27).Value = 0.002335768728363975
$ws.Cells.Item(6, 28).Value = 0.002106293295639031
$ws.Cells.Item(6, 29).Value = 0.0002615803796677384
$ws.Cells.Item(6, 30).Value = 0.09582973492160531
$ws.Cells.Item(6, 31).Value = 0.1401594763039248
$ws.Cells.Item(6, 32).Value = 0.1220543389049664
